$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3
